{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\n\n// Locate the \"Diagrama de Colaboraci\u00f3n\" section heading (PSI-T\u00edtulo2) that\n// starts the block to remove, and the final \"PSI-T\u00edtulo1\" paragraph whose\n// picture must be cleared but whose (now empty) paragraph must remain.\nparas.load(\"items/text,items/style\");\nawait context.sync();\n\nconst items = paras.items;\nlet startIndex = -1; // first paragraph of the block to delete\nlet endIndex = -1;   // last paragraph of the block to delete (inclusive)\n\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].style === \"PSI - T\u00edtulo 2\" && items[i].text === \"Diagrama de Colaboraci\u00f3n\") {\n    // The block actually begins at the blank \"PSI - T\u00edtulo 2\" paragraph\n    // immediately preceding this heading (it only carries the bookmark).\n    startIndex = i;\n    if (i - 1 >= 0 && items[i - 1].style === \"PSI - T\u00edtulo 2\" && items[i - 1].text === \"\") {\n      startIndex = i - 1;\n    }\n    break;\n  }\n}\n\nfor (let i = startIndex; i < items.length; i++) {\n  if (items[i].style === \"PSI - T\u00edtulo 1\") {\n    endIndex = i - 1;\n    break;\n  }\n}\n\nif (startIndex !== -1 && endIndex !== -1 && endIndex >= startIndex) {\n  for (let i = endIndex; i >= startIndex; i--) {\n    items[i].delete();\n  }\n}\nawait context.sync();\n\n// Remove the picture that used to sit in the trailing \"PSI - T\u00edtulo 1\"\n// paragraph, leaving the (now empty) paragraph itself untouched.\nconst paras2 = body.paragraphs;\nparas2.load(\"items/style\");\nawait context.sync();\n\nconst items2 = paras2.items;\nfor (let i = items2.length - 1; i >= 0; i--) {\n  if (items2[i].style === \"PSI - T\u00edtulo 1\") {\n    const pics = items2[i].inlinePictures;\n    pics.load(\"items\");\n    await context.sync();\n    for (let j = pics.items.length - 1; j >= 0; j--) {\n      pics.items[j].delete();\n    }\n    await context.sync();\n    break;\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the block to delete: it begins at the blank \"PSI - T\u00edtulo 2\"\n# paragraph that precedes the \"Diagrama de Colaboraci\u00f3n\" heading (that\n# blank paragraph only carries the section bookmark) and runs up to, but\n# not including, the trailing \"PSI - T\u00edtulo 1\" paragraph.\n$startIndex = -1\n$endIndex = -1\n$n = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $n; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $style = $p.Style.NameLocal\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($style -eq \"PSI - T\u00edtulo 2\" -and $text -eq \"Diagrama de Colaboraci\u00f3n\") {\n        $startIndex = $i\n        if ($i -gt 1) {\n            $prev = $d.Paragraphs.Item($i - 1)\n            $prevText = $prev.Range.Text.TrimEnd([char]13, [char]7)\n            if ($prev.Style.NameLocal -eq \"PSI - T\u00edtulo 2\" -and $prevText -eq \"\") {\n                $startIndex = $i - 1\n            }\n        }\n        break\n    }\n}\n\nif ($startIndex -ne -1) {\n    for ($i = $startIndex; $i -le $n; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Style.NameLocal -eq \"PSI - T\u00edtulo 1\") {\n            $endIndex = $i - 1\n            break\n        }\n    }\n}\n\nif ($startIndex -ne -1 -and $endIndex -ne -1 -and $endIndex -ge $startIndex) {\n    $startPara = $d.Paragraphs.Item($startIndex)\n    $endPara = $d.Paragraphs.Item($endIndex)\n    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $range.Delete()\n}\n\n# Remove the picture that used to sit in the trailing \"PSI - T\u00edtulo 1\"\n# paragraph, leaving the (now empty) paragraph itself in place. Deleting\n# through a Range's own InlineShapes collection is unreliable here, so the\n# matching shape is located by position and removed via the document-level\n# InlineShapes collection instead.\n$n2 = $d.Paragraphs.Count\n$targetPara = $null\nfor ($i = $n2; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Style.NameLocal -eq \"PSI - T\u00edtulo 1\") {\n        $targetPara = $p\n        break\n    }\n}\n\nif ($targetPara -ne $null) {\n    $pStart = $targetPara.Range.Start\n    $pEnd = $targetPara.Range.End\n    $allShapes = $d.InlineShapes\n    for ($k = $allShapes.Count; $k -ge 1; $k--) {\n        $shape = $allShapes.Item($k)\n        if ($shape.Range.Start -ge $pStart -and $shape.Range.End -le $pEnd) {\n            $shape.Delete()\n        }\n    }\n}\n"}
